$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(0.2225774455180299, 0.6025194276978968, 0.9375366183193465, 0.968264745986007, 0.9517120452360237, 51)
    3  = @(0.5002927602626792, 0.8093331159883599, 2.512481039783275, 1.585080767589865, 1.519327252374612, 50)
    4  = @(0.7210151913624349, 1.262988772388277, 6.435165304488891, 2.536762760781719, 2.457343660326128, 49)
    5  = @(0.4586541891000874, 1.384892259048143, 7.994179379982225, 2.82739798754654, 2.819473144954236, 48)
    6  = @(0.4309535499285895, 1.407201242206765, 8.132525087741485, 2.85175824496774, 2.84948429316084, 47)
    7  = @(0.436267056338538, 1.569053982142953, 9.223627919858613, 3.037042627270584, 3.045889392437618, 38)
    8  = @(0.4110347441899493, 1.588476331609555, 9.491981118609186, 3.080905892527259, 3.095481328727775, 37)
    9  = @(0.08076273096934439, 2.109961141030493, 15.00836014869273, 3.874062486420776, 3.973840450533185, 20)
    10 = @(-0.6203614576892605, 2.00672458876966, 12.86991326340611, 3.587466134112782, 3.677701320433923, 13)
    11 = @(1.509871661924674, 1.524567153826987, 3.867561475191954, 1.966611673714959, 1.408833311515586, 5)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Cells.Item($row, 2).Value = $rowValues[0]
    $ws.Cells.Item($row, 3).Value = $rowValues[1]
    $ws.Cells.Item($row, 4).Value = $rowValues[2]
    $ws.Cells.Item($row, 5).Value = $rowValues[3]
    $ws.Cells.Item($row, 6).Value = $rowValues[4]
    $ws.Cells.Item($row, 7).Value = $rowValues[5]
}

$wb.Save()
